$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)  # yozgat

$ws.Cells.Item(3, 4).Value = 0.01
$ws.Cells.Item(3, 5).Value = 0.25
$ws.Cells.Item(3, 6).Value = 47.3
$ws.Cells.Item(3, 7).Value = 4.88
$ws.Cells.Item(4, 4).Value = 0.04
$ws.Cells.Item(4, 5).Value = 0.73
$ws.Cells.Item(4, 6).Value = 12.79
$ws.Cells.Item(4, 7).Value = 1.76
$ws.Cells.Item(5, 4).Value = 0.05
$ws.Cells.Item(5, 5).Value = 1.04
$ws.Cells.Item(5, 6).Value = 35.42
$ws.Cells.Item(5, 7).Value = 12.68
$ws.Cells.Item(6, 4).Value = 0.04
$ws.Cells.Item(6, 5).Value = 0.89
$ws.Cells.Item(6, 6).Value = 32.04
$ws.Cells.Item(6, 7).Value = 4.2
$ws.Cells.Item(7, 4).Value = 0.05
$ws.Cells.Item(7, 5).Value = 1.01
$ws.Cells.Item(7, 6).Value = 8.46
$ws.Cells.Item(7, 7).Value = 1.69
$ws.Cells.Item(8, 4).Value = 0.04
$ws.Cells.Item(8, 5).Value = 0.85
$ws.Cells.Item(8, 6).Value = 55.84
$ws.Cells.Item(8, 7).Value = 2.91
$ws.Cells.Item(9, 4).Value = 0.05
$ws.Cells.Item(9, 5).Value = 1.04
$ws.Cells.Item(9, 6).Value = 112.61
$ws.Cells.Item(9, 7).Value = 7.86
$ws.Cells.Item(10, 4).Value = 0.03
$ws.Cells.Item(10, 5).Value = 0.62
$ws.Cells.Item(10, 6).Value = 37.05
$ws.Cells.Item(10, 7).Value = 21.5
$ws.Cells.Item(11, 4).Value = 0.03
$ws.Cells.Item(11, 5).Value = 0.59
$ws.Cells.Item(11, 6).Value = 65.05
$ws.Cells.Item(11, 7).Value = 9.42
$ws.Cells.Item(12, 4).Value = 0.08
$ws.Cells.Item(12, 5).Value = 1.55
$ws.Cells.Item(12, 6).Value = 39.2
$ws.Cells.Item(12, 7).Value = 5.69
$ws.Cells.Item(13, 4).Value = 0.05
$ws.Cells.Item(13, 5).Value = 1.04
$ws.Cells.Item(13, 6).Value = 81.66
$ws.Cells.Item(13, 7).Value = 2.98
$ws.Cells.Item(14, 4).Value = 0.01
$ws.Cells.Item(14, 5).Value = 0.17
$ws.Cells.Item(14, 6).Value = 37.76
$ws.Cells.Item(14, 7).Value = 4.61
$ws.Cells.Item(15, 4).Value = 0.06
$ws.Cells.Item(15, 5).Value = 1.25
$ws.Cells.Item(15, 6).Value = 32.48
$ws.Cells.Item(15, 7).Value = 2.98
$ws.Cells.Item(16, 4).Value = 0.08
$ws.Cells.Item(16, 5).Value = 1.6
$ws.Cells.Item(16, 6).Value = 166.88
$ws.Cells.Item(16, 7).Value = 4.61
$ws.Cells.Item(17, 4).Value = 0.08
$ws.Cells.Item(17, 5).Value = 1.54
$ws.Cells.Item(17, 6).Value = 203.81
$ws.Cells.Item(17, 7).Value = 5.76
$ws.Cells.Item(18, 4).Value = 0.06
$ws.Cells.Item(18, 5).Value = 1.2
$ws.Cells.Item(18, 6).Value = 34.05
$ws.Cells.Item(18, 7).Value = 13.02
$ws.Cells.Item(19, 4).Value = 0.05
$ws.Cells.Item(19, 5).Value = 1.04
$ws.Cells.Item(19, 6).Value = 17.67
$ws.Cells.Item(19, 7).Value = 20.48
$ws.Cells.Item(20, 4).Value = 0.04
$ws.Cells.Item(20, 5).Value = 0.79
$ws.Cells.Item(20, 6).Value = 71.67
$ws.Cells.Item(20, 7).Value = 3.93
$ws.Cells.Item(21, 4).Value = 0.07
$ws.Cells.Item(21, 5).Value = 1.4
$ws.Cells.Item(21, 6).Value = 35.42
$ws.Cells.Item(21, 7).Value = 1.48
$ws.Cells.Item(22, 4).Value = 0.08
$ws.Cells.Item(22, 5).Value = 1.55
$ws.Cells.Item(22, 6).Value = 73.94
$ws.Cells.Item(22, 7).Value = 1.82
$ws.Cells.Item(23, 4).Value = 0.04
$ws.Cells.Item(23, 5).Value = 0.83
$ws.Cells.Item(23, 6).Value = 52.73
$ws.Cells.Item(23, 7).Value = 2.71
$ws.Cells.Item(24, 4).Value = 0.06
$ws.Cells.Item(24, 5).Value = 1.24
$ws.Cells.Item(24, 6).Value = 57.28
$ws.Cells.Item(24, 7).Value = 6.78
$ws.Cells.Item(25, 4).Value = 0.07
$ws.Cells.Item(25, 5).Value = 1.44
$ws.Cells.Item(25, 6).Value = 69.11
$ws.Cells.Item(25, 7).Value = 3.66
$ws.Cells.Item(26, 4).Value = 0.03
$ws.Cells.Item(26, 5).Value = 0.51
$ws.Cells.Item(26, 6).Value = 10.05
$ws.Cells.Item(26, 7).Value = 2.5
$ws.Cells.Item(27, 4).Value = 0.01
$ws.Cells.Item(27, 5).Value = 0.19
$ws.Cells.Item(27, 6).Value = 16.02
$ws.Cells.Item(27, 7).Value = 1.28

$ws.Activate()
$ws.Range("D3:G27").Select()

$ws1 = $wb.Worksheets.Item(1)  # adana
$ws1.Activate()
